# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, styled like the other header cells (B1:G1) -
# copy G1's formatting (bordered/bold/centered header style) into H1, then set its text
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Values for H2:H8 (the "Save" stat per row)
$values = @(0, 0, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
